$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignInTest")
$ws.Range("C2").Value = "Y"
